# Slide 13, shape 3 ("REFERENCE" body placeholder) holds a numbered list of
# citations. The 5th citation (Jiuxiang Gu et al., Stack-Captioning) currently
# ends with a hyperlinked URL run; replace that hyperlink run with the plain
# text "Image Captioning" and drop the trailing "Image Captioning " wording
# that was duplicated in the first (citation) run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)
$sh = $s.Shapes.Item(3)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para = $tr.Paragraphs(5)

$newText = "Jiuxiang Gu, Jianfei Cai, Gang Wang and Tsuhan Chen [2018]. Stack-Captioning: Coarse-to-Fine Learning for Image Captioning"

# Collapse the whole paragraph down to a single run carrying the first run's
# (plain, non-hyperlinked) formatting -- this is what drops the old second
# run's underline / hyperlink-color / hlinkClick.
$para.Text = $newText

# Re-split off the trailing "Image Captioning" into its own run (still with
# the same plain formatting) so it matches the two-run structure of the
# original citation (citation text run + trailing run).
$tailStart = $newText.Length - "Image Captioning".Length + 1
$tail = $para.Characters($tailStart, "Image Captioning".Length)
$tail.Text = "Image Captioning"
